$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("B2").Value = 3900
$ws.Range("B3").Value = 4950

# Update the active selection to B3 (as reflected in the saved view state)
$ws.Range("B3").Select()
